# Apply latest cryptocurrency price/volume snapshot (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.670.98"
$ws.Range("E2").Value = "  +3.76%  "
$ws.Range("D3").Value = "1.913.01"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.01"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.700"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.35"
$ws.Range("E8").Value = "  +4.25%  "
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.66"
$ws.Range("E10").Value = "  +9.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0762"
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.56"
$ws.Range("E13").Value = "  +7.37%  "
$ws.Range("E14").Value = "  +6.87%  "
$ws.Range("D15").Value = "2.191.66"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("D17").Value = "1.919.46"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "36.630.65"
$ws.Range("E18").Value = "  +3.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.77"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "0.0₃0861"
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "250.23"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.37"
$ws.Range("E22").Value = "  +4.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.63"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.19"
$ws.Range("E27").Value = "  +2.20%  "
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.70"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.56"
$ws.Range("E31").Value = "  +6.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0618"
$ws.Range("E32").Value = "  +3.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.33"
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0892"
$ws.Range("E34").Value = "  +22.68%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.876"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.76"
$ws.Range("E39").Value = "  +49.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.01"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.32"
$ws.Range("E41").Value = "  +8.41%  "
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.40"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.92"
$ws.Range("E44").Value = "  +21.91%  "
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").Value = "1.345.22"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.45"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.38"
$ws.Range("E51").Value = "  +2.94%  "

Write-Output "Updated 89 cells"
